# Apply forecast data updates (Optuna attempt - go back with original)
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet updates ---
# Row 2 (W8)
$wsForecast.Range("D2").Value = 153
$wsForecast.Range("H2").Value = 10.75
$wsForecast.Range("L2").Value = 1.13

# Row 3 (W9)
$wsForecast.Range("D3").Value = 150
$wsForecast.Range("H3").Value = 9.94
$wsForecast.Range("L3").Value = 0.88

# Row 4 (W10)
$wsForecast.Range("H4").Value = 8.890000000000001
$wsForecast.Range("L4").Value = 1.17

# Row 5 (W11)
$wsForecast.Range("H5").Value = 7.76
$wsForecast.Range("L5").Value = 0.92

# Row 6 (W12)
$wsForecast.Range("H6").Value = 6.7
$wsForecast.Range("L6").Value = 1.06

# Row 7 (W13)
$wsForecast.Range("H7").Value = 5.8
$wsForecast.Range("L7").Value = 0.9399999999999999

# Row 8 (W14)
$wsForecast.Range("D8").Value = 147
$wsForecast.Range("H8").Value = 4.97
$wsForecast.Range("L8").Value = 1.17

# Row 9 (W15)
$wsForecast.Range("D9").Value = 131
$wsForecast.Range("H9").Value = 4.45
$wsForecast.Range("L9").Value = 0.82

# Row 10 (W16)
$wsForecast.Range("D10").Value = 133
$wsForecast.Range("H10").Value = 3.4
$wsForecast.Range("L10").Value = 0.97

# Row 11 (W17)
$wsForecast.Range("D11").Value = 142
$wsForecast.Range("H11").Value = 2.25
$wsForecast.Range("L11").Value = 0.89

# Row 12 (W18)
$wsForecast.Range("D12").Value = 140
$wsForecast.Range("H12").Value = 1.26
$wsForecast.Range("J12").Value = "Normal"
$wsForecast.Range("L12").Value = 1.2

# Row 13 (W19)
$wsForecast.Range("D13").Value = 116
$wsForecast.Range("H13").Value = 0.32
$wsForecast.Range("L13").Value = 0.9

# Row 14 (W20)
$wsForecast.Range("D14").Value = 91
$wsForecast.Range("L14").Value = 0.92

# Row 15 (W21)
$wsForecast.Range("D15").Value = 84
$wsForecast.Range("L15").Value = 1.05

# Row 16 (W22)
$wsForecast.Range("D16").Value = 99
$wsForecast.Range("L16").Value = 0.87

# Row 17 (W23)
$wsForecast.Range("D17").Value = 117
$wsForecast.Range("L17").Value = 0.96

# --- Summary sheet updates ---
# Keep these cells as text (matching original inlineStr type) instead of
# letting Excel auto-convert numeric-looking strings to numbers.
$wsSummary.Range("B9:B12").NumberFormat = "@"
$wsSummary.Range("B14").NumberFormat = "@"

$wsSummary.Range("B9").Value = "2114"
$wsSummary.Range("B10").Value = "1192"
$wsSummary.Range("B11").Value = "607"
$wsSummary.Range("B12").Value = "155"
$wsSummary.Range("B14").Value = "84"
